$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# --- Row 16 (Precinct row "23") ---
# C16/D16/E16 switch from numeric cells to the "N/A" text style (s=14) used
# elsewhere in the sheet (e.g. C22/D22/E22), so copy style+value from those
# donor cells instead of just assigning a literal value (which Excel would
# otherwise coerce back to a number).
$ws.Range("C22").Copy($ws.Range("C16"))
$ws.Range("D22").Copy($ws.Range("D16"))
$ws.Range("E22").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = 0
$ws.Range("N16").Value = -66.666666666666

# --- Row 17 (Precinct row "24") ---
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 60
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 20
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 71.428571428571

# --- Row 18 (Precinct row "25") ---
# D18 and E18 switch from the "N/A" text style back to numeric styles, so
# copy style+value from a same-column/compatible numeric donor cell first,
# then overwrite with the real number.
$ws.Range("D17").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("N14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -60
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = -72.727272727272
$ws.Range("M18").Value = -83.333333333333
$ws.Range("N18").Value = -93.181818181818

# --- Row 19 (Precinct row "26") ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 40
$ws.Range("J19").Value = 43
$ws.Range("K19").Value = -6.976744186046
$ws.Range("L19").Value = -11.111111111111
$ws.Range("M19").Value = 53.846153846153
$ws.Range("N19").Value = 90.476190476190

# --- Row 20 (Precinct row "27") ---
# C20 switches from the "N/A" text style to a numeric style; use F20 (same
# numeric style, unrelated column) as the style donor before overwriting it
# with the real number.
$ws.Range("F20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 8
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = -57.894736842105
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -92.727272727272

# --- Row 21 (Precinct row "28", TOTAL row - bold style) ---
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 38
$ws.Range("H21").Value = 10.526315789473
$ws.Range("I21").Value = 69
$ws.Range("J21").Value = 77
$ws.Range("K21").Value = -10.389610389610
$ws.Range("L21").Value = -8
$ws.Range("M21").Value = 9.523809523809
$ws.Range("N21").Value = -64.0625

# --- Row 24 (Precinct row "31", Petit Larceny) ---
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 25
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = -19.444444444444
$ws.Range("I24").Value = 58
$ws.Range("J24").Value = 81
$ws.Range("K24").Value = -28.395061728395
$ws.Range("L24").Value = 7.407407407407
$ws.Range("M24").Value = -28.395061728395

# --- Row 25 (Precinct row "32", Misd. Assault) ---
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -66.666666666666
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -38.095238095238
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 33
$ws.Range("K25").Value = -39.393939393939
$ws.Range("L25").Value = -23.076923076923
$ws.Range("M25").Value = -35.483870967741
